$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (existing "VA.MHV.PHR.vitals" profile row) is narrowed: Time Types
# drops "Periodĵ", and Data Absent Reason becomes "prohibited" now that
# BP/Pain have their own rows below.
$ws.Range("G3").Value = "dateTimeĵ"
$ws.Range("I3").Value = "prohibited"

# Stamp the existing row-3 cell format (border/fill/wrap) onto the four new
# rows first, so every new cell (including the ones left blank) carries the
# same style index as the rest of the table.
$ws.Range("A3:K3").Copy()
$ws.Range("A4:K7").PasteSpecial(-4122)

# New row 4: Blood Pressure panel profile
$ws.Range("A4").Value = "VA.MHV.PHR.vitalsBP"
$ws.Range("B4").Value = "VA MHV PHR Vital-Signs for Blood Pressure"
$ws.Range("C4").Value = "Observation Category Codes#vital-signs"
$ws.Range("E4").Value = "LOINC#85354-9"
$ws.Range("G4").Value = "dateTimeĵ, Periodĵ"
$ws.Range("I4").Value = "optional"

# New row 5: Blood Pressure - systolic component
$ws.Range("B5").Value = "VA MHV PHR Vital-Signs for Blood Pressure"
$ws.Range("E5").Value = "LOINC#8480-6"
$ws.Range("H5").Value = "Quantityĵ"
$ws.Range("I5").Value = "optional"

# New row 6: Blood Pressure - diastolic component
$ws.Range("B6").Value = "VA MHV PHR Vital-Signs for Blood Pressure"
$ws.Range("E6").Value = "LOINC#8462-4"
$ws.Range("H6").Value = "Quantityĵ"
$ws.Range("I6").Value = "optional"

# New row 7: Pain vital-signs profile
$ws.Range("A7").Value = "VA.MHV.PHR.vitalsPain"
$ws.Range("B7").Value = "VA MHV PHR Vital-Signs for PAIN"
$ws.Range("C7").Value = "Observation Category Codes#vital-signs"
$ws.Range("F7").Value = "http://hl7.org/fhir/us/core/ValueSet/us-core-vital-signs (extensible)"
$ws.Range("G7").Value = "dateTimeĵ, Periodĵ"
$ws.Range("H7").Value = "Quantityĵ"
$ws.Range("I7").Value = "optional"
